$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add new Sheet2 right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Fill in Sheet2 data
$data = @(
    @(1, "P", 2, 4, 5, 3),
    @(1, "Q", 7, 3, 10, 2),
    @(3, "P", 4, 10, 15, 1),
    @(3, "Q", 6, 12, 8, 12)
)

for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Update selections
[void]$ws1.Range("A1:D4").Select()
[void]$ws2.Range("C6").Select()

# Make Sheet2 the active sheet/tab
[void]$ws2.Activate()
